# ============================================================
# Apply "se agrego todo lo de tips" commit:
#  - add "tips" sheet (destinos tips/consejos de viaje)
#  - add "iconos" sheet (icon name -> FontAwesome class lookup)
#  - update selection on "destinos" and make "iconos" the active tab
# ============================================================

$wb = $excel.ActiveWorkbook

# --- 1) "destinos": just move the selection, keep all data/formatting as-is ---
$ws1 = $wb.Worksheets.Item("destinos")
$ws1.Range("B14").Select()

# --- 2) Add the "tips" sheet right after "destinos" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws2.Name = "tips"

# --- sheet2 "tips" data ---
$ws2.Range("A1").Value = "titulo"
$ws2.Range("B1").Value = "subtitulo"
$ws2.Range("C1").Value = "descripcion"
$ws2.Range("D1").Value = "icono_id"
$ws2.Range("A2").Value = "Transporte Público"
$ws2.Range("B2").Value = "¡Tu huella también cuenta!"
$ws2.Range("C2").Value = "Al elegir transporte público, como colectivos, subtes o trenes, en vez de vehículos particulares, reducís significativamente tus emisiones de carbono.`nPlanificar tu viaje con conciencia ambiental también es parte del compromiso con el planeta."
$ws2.Range("A3").Value = "Optá por productos reutilizables."
$ws2.Range("B3").Value = "Optá por productos reutilizables."
$ws2.Range("C3").Value = "Llevá siempre contigo una botella de agua reutilizable, un termo para el café y bolsas de tela. Cada acción suma para evitar el uso de plásticos.`nSi no podés evitar el plástico, trata de reciclarlo o llevarlo a puntos de recolección selectiva."
$ws2.Range("A4").Value = "Movilidad Activa"
$ws2.Range("B4").Value = "¡Usá la bici y recorre la ciudad!"
$ws2.Range("C4").Value = "Camina o usa bicicleta para distancias cortas.`nLa Ciudad de Buenos Aires tiene un sistema de bicicletas públicas gratuitas llamado EcoBici, que te permite recorrer el centro y otras zonas de forma gratuita.`n¡Solo necesitás registrarte en la app y comenzar a pedalear!"
$ws2.Range("A5").Value = "Documentación Digital"
$ws2.Range("C5").Value = "Llevá pasajes, reservas y documentos en el celular en lugar de imprimirlos.`nHoy en día la mayoría de aerolíneas, hoteles y transportes aceptan comprobantes digitales. Además de ahorrar papel, evitás perder documentación durante el viaje."
$ws2.Range("A6").Value = "Hidratación consciente"
$ws2.Range("B6").Value = "Tomá agua sin desperdiciar"
$ws2.Range("C6").Value = "Llevá siempre una botella reutilizable y evitá comprar botellas descartables.`nEn muchas zonas de Buenos Aires podés recargar agua potable en bares, estaciones y espacios públicos.`nAdemás de cuidarte, ayudás a reducir el consumo de plástico y el gasto innecesario de agua."
$ws2.Range("A7").Value = "`nApoya lo Local"
$ws2.Range("B7").Value = "¡Usá la bici y recorre la ciudad!"
$ws2.Range("C7").Value = "Compra en mercados locales y apoya proyectos comunitarios sustentables"
$ws2.Range("A8").Value = "Cuidado con el consumo masivo"
$ws2.Range("B8").Value = "Comprá solo lo necesario."
$ws2.Range("C8").Value = "Evitá comprar productos descartables o souvenirs de uso innecesario.`nPriorizar calidad sobre cantidad reduce residuos y el impacto ambiental del turismo masivo.`nElegir con conciencia también es una forma de cuidar el planeta."
$ws2.Range("A9").Value = "Respetá la naturaleza"
$ws2.Range("B9").Value = "Respetá cada lugar que visitás"
$ws2.Range("C9").Value = "No dejes residuos, respetá la flora y la fauna y mantenete en los senderos habilitados.`nCuidar los espacios naturales permite que más personas puedan disfrutarlos hoy y en el futuro.`nViajar también es aprender a convivir con el entorno."
$ws2.Range("A10").Value = "Cuida el Agua"
$ws2.Range("B10").Value = "Menos desperdicio de agua"
$ws2.Range("C10").Value = "Evitá dejar correr el agua innecesariamente y cerrá bien las canillas luego de usarlas.`nPequeñas acciones, como no desperdiciar agua al lavar objetos o manos,ayudan a conservar uno de los recursos más importantes del planeta."
$ws2.Range("A11").Value = "Gestioná tus residuos"
$ws2.Range("B11").Value = "¡Reducí, reutilizá y recicla!"
$ws2.Range("C11").Value = "Separá residuos siempre que sea posible y evitá tirar basura en espacios públicos.`nMuchos destinos cuentan con puntos de reciclaje que ayudan a reducir el impacto ambiental.`nCada residuo bien gestionado suma."
$ws2.Range("A12").Value = "Ahorrá energía"
$ws2.Range("B12").Value = "Pequeñas acciones, gran impacto"
$ws2.Range("C12").Value = "Apagá luces, aire acondicionado y dispositivos cuando no los estés usando.`nReducir el consumo de energía ayuda a disminuir las emisiones y el gasto innecesario de recursos.`nSer responsable también es parte del viaje."

# --- sheet2 row heights + wrap text on column C (long descriptions) ---
$ws2.Range("C2").WrapText = $true
$ws2.Rows.Item(2).RowHeight = 30
$ws2.Range("C3").WrapText = $true
$ws2.Rows.Item(3).RowHeight = 30
$ws2.Range("C4").WrapText = $true
$ws2.Rows.Item(4).RowHeight = 60
$ws2.Range("C5").WrapText = $true
$ws2.Rows.Item(5).RowHeight = 45
$ws2.Range("C6").WrapText = $true
$ws2.Rows.Item(6).RowHeight = 45
$ws2.Range("C7").WrapText = $true
$ws2.Rows.Item(7).RowHeight = 30
$ws2.Range("C8").WrapText = $true
$ws2.Rows.Item(8).RowHeight = 45
$ws2.Range("C9").WrapText = $true
$ws2.Rows.Item(9).RowHeight = 45
$ws2.Range("C10").WrapText = $true
$ws2.Rows.Item(10).RowHeight = 30
$ws2.Range("C11").WrapText = $true
$ws2.Rows.Item(11).RowHeight = 45
$ws2.Range("C12").WrapText = $true
$ws2.Rows.Item(12).RowHeight = 45

# --- sheet2 column widths ---
$ws2.Range("A:B").ColumnWidth = 30.43
$ws2.Range("C:C").ColumnWidth = 135

$ws2.Range("C13").Select()

# --- 3) Add the "iconos" sheet right after "tips" ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$ws3.Name = "iconos"

# --- sheet3 "iconos" data ---
$ws3.Range("A1").Value = "nombre"
$ws3.Range("B1").Value = "clase"
$ws3.Range("A2").Value = "Reciclaje"
$ws3.Range("B2").Value = "fa-solid fa-recycle"
$ws3.Range("A3").Value = "Naturaleza"
$ws3.Range("B3").Value = "fa-solid fa-seedling"
$ws3.Range("A4").Value = "Digital"
$ws3.Range("B4").Value = "fa-solid fa-mobile-screen"
$ws3.Range("A5").Value = "Transporte"
$ws3.Range("B5").Value = "fa-solid fa-bus"
$ws3.Range("A6").Value = "Residuos"
$ws3.Range("B6").Value = "fa-solid fa-trash-can-arrow-up"
$ws3.Range("A7").Value = "Energía"
$ws3.Range("B7").Value = "fa-regular fa-lightbulb"
$ws3.Range("A8").Value = "Mapa"
$ws3.Range("B8").Value = "fa-solid fa-map"
$ws3.Range("A9").Value = "Hotel"
$ws3.Range("B9").Value = "fa-solid fa-hotel"
$ws3.Range("A10").Value = "Bicicleta"
$ws3.Range("B10").Value = "fa-solid fa-person-biking"
$ws3.Range("A11").Value = "Caminata"
$ws3.Range("B11").Value = "fa-solid fa-person-walking"

# --- sheet3 column widths ---
$ws3.Range("A:A").ColumnWidth = 15.43
$ws3.Range("B:B").ColumnWidth = 35.71

# "iconos" ends up the active tab / selected sheet, matching the target workbook view
$ws3.Range("A2").Select()
